$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    # Forces the cell to store $val as literal text even when it
    # looks numeric (e.g. "586.43"), matching the source inlineStr cells,
    # then restores the default (unstyled) cell style so no stray
    # number-format style is left behind on the cell.
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "63.490.26"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").Value = "3.063.43"
$ws.Range("E3").Value = "  -4.08%  "
$ws.Range("E4").Value = "  -0.16%  "
Set-TextValue "D5" "586.43"
$ws.Range("E5").Value = "  -1.37%  "
Set-TextValue "D6" "153.55"
$ws.Range("E6").Value = "  +3.02%  "
$ws.Range("E7").Value = "  -0.10%  "
Set-TextValue "D8" "0.532"
$ws.Range("E8").Value = "  -0.47%  "
$ws.Range("D9").Value = "3.063.57"
$ws.Range("E9").Value = "  -3.77%  "
Set-TextValue "D10" "0.155"
$ws.Range("E10").Value = "  -5.37%  "
$ws.Range("E11").Value = "  -2.98%  "
Set-TextValue "D12" "0.447"
$ws.Range("E12").Value = "  -3.61%  "
Set-TextValue "D13" "36.61"
$ws.Range("E13").Value = "  -3.45%  "
Set-TextValue "D14" "0.0000236"
$ws.Range("E14").Value = "  -5.42%  "
$ws.Range("E15").Value = "  -2.48%  "
$ws.Range("D16").Value = "3.567.14"
$ws.Range("E16").Value = "  -4.09%  "
$ws.Range("D17").Value = "63.468.78"
$ws.Range("E17").Value = "  -1.23%  "
Set-TextValue "D18" "7.09"
$ws.Range("E18").Value = "  -3.91%  "
$ws.Range("D19").Value = "3.064.32"
$ws.Range("E19").Value = "  -3.99%  "
Set-TextValue "D20" "467.48"
$ws.Range("E20").Value = "  -1.43%  "
Set-TextValue "D21" "14.22"
$ws.Range("E21").Value = "  -2.87%  "
Set-TextValue "D22" "0.699"
$ws.Range("E22").Value = "  -5.65%  "
$ws.Range("E23").Value = "  -3.48%  "
Set-TextValue "D24" "2.41"
$ws.Range("E24").Value = "  -0.99%  "
Set-TextValue "D25" "80.15"
$ws.Range("E25").Value = "  -2.19%  "
Set-TextValue "D26" "12.68"
$ws.Range("E26").Value = "  -4.15%  "
Set-TextValue "D27" "10.37"
$ws.Range("E27").Value = "  +3.51%  "
$ws.Range("E28").Value = "  -0.27%  "
Set-TextValue "D29" "7.31"
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("E30").Value = "  -0.20%  "
Set-TextValue "D31" "2.64"
$ws.Range("E31").Value = "  -3.78%  "
Set-TextValue "D32" "2.13"
$ws.Range("E32").Value = "  -6.18%  "
Set-TextValue "D33" "26.94"
$ws.Range("E33").Value = "  -5.71%  "
Set-TextValue "D34" "0.110"
$ws.Range("E34").Value = "  -7.85%  "
$ws.Range("D35").Value = "0.0₃0815"
$ws.Range("E35").Value = "  -6.04%  "
$ws.Range("E36").Value = "  -2.84%  "
Set-TextValue "D37" "5.94"
$ws.Range("E37").Value = "  -5.00%  "
Set-TextValue "D38" "3.23"
$ws.Range("E38").Value = "  -4.93%  "
Set-TextValue "D39" "2.20"
$ws.Range("E39").Value = "  -5.58%  "
Set-TextValue "D40" "50.43"
$ws.Range("E40").Value = "  -2.22%  "
Set-TextValue "D41" "9.13"
$ws.Range("E41").Value = "  -3.43%  "
Set-TextValue "D42" "435.11"
$ws.Range("E42").Value = "  -7.34%  "
Set-TextValue "D43" "0.285"
$ws.Range("E43").Value = "  -3.80%  "
Set-TextValue "D44" "40.24"
$ws.Range("E44").Value = "  +1.87%  "
$ws.Range("E45").Value = "  +1.16%  "
Set-TextValue "D46" "0.0356"
$ws.Range("E46").Value = "  -6.09%  "
$ws.Range("D47").Value = "2.789.61"
$ws.Range("E47").Value = "  -4.77%  "
Set-TextValue "D48" "128.83"
$ws.Range("E48").Value = "  -3.37%  "
Set-TextValue "D49" "1.00"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("E50").Value = "  +0.90%  "
Set-TextValue "D51" "2.20"
$ws.Range("E51").Value = "  -3.53%  "
